$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.242356181144714
$ws.Range("B1").Value = 2.098983526229858
$ws.Range("C1").Value = 5.932182312011719
$ws.Range("D1").Value = 1.966402649879456
$ws.Range("E1").Value = 1.142043113708496
